# Added SVR parameter loading from pred_par structure and Excel files:
# extend the pred_par table with three new columns (K:M) holding the
# SVR hyper-parameters, mirroring the layout of the existing
# rnn / data columns in row 1 (headers) and row 2 (values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers (row 1)
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# New parameter values (row 2)
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# Reflect the new active cell/selection left after entering the data
$ws.Range("K7").Select() | Out-Null
